$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-07 Saturday" "2026-02-08 Sunday"

Replace-Text "975÷7=139, 2" "317÷8=39, 5"
Replace-Text "743÷9=82, 5" "269÷2=134, 1"
Replace-Text "829÷6=138, 1" "710÷5=142, 0"
Replace-Text "776÷3=258, 2" "253÷9=28, 1"
Replace-Text "520÷7=74, 2" "886÷3=295, 1"

Replace-Text "239÷6=39, 5" "815÷5=163, 0"
Replace-Text "721÷4=180, 1" "513÷8=64, 1"
Replace-Text "985÷5=197, 0" "707÷9=78, 5"
Replace-Text "715÷4=178, 3" "841÷8=105, 1"
Replace-Text "293÷7=41, 6" "603÷2=301, 1"

Replace-Text "856÷6=142, 4" "156÷9=17, 3"
Replace-Text "471÷8=58, 7" "929÷3=309, 2"
Replace-Text "367÷4=91, 3" "809÷3=269, 2"
Replace-Text "280÷9=31, 1" "231÷3=77, 0"
Replace-Text "891÷6=148, 3" "468÷6=78, 0"

Replace-Text "197÷4=49, 1" "101÷2=50, 1"
Replace-Text "232÷2=116, 0" "517÷6=86, 1"
Replace-Text "926÷3=308, 2" "465÷4=116, 1"
Replace-Text "546÷2=273, 0" "648÷3=216, 0"
Replace-Text "784÷9=87, 1" "530÷8=66, 2"

Replace-Text "487÷7=69, 4" "304÷8=38, 0"
Replace-Text "641÷9=71, 2" "765÷8=95, 5"
Replace-Text "692÷2=346, 0" "994÷9=110, 4"
Replace-Text "451÷4=112, 3" "467÷4=116, 3"
Replace-Text "133÷2=66, 1" "786÷4=196, 2"
